$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-6 from 45183 (2023-09-14)
# to 45184 (2023-09-15), keeping the existing date formatting/style intact.
foreach ($r in 2..6) {
    $ws.Cells.Item($r, 3).Value = 45184
}
